$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the wind_excess_energy_price / solar_excess_energy_price columns (old E:F),
# shifting all subsequent columns left by two.
$ws.Range("E:F").Delete()

# Restore selection to B3 (as captured in the saved workbook state).
$ws.Range("B3").Select()
